$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.126.55"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "3.328.84"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.08"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.11"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.66%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.323.11"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.575"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.62%  "
$ws.Range("E10").Value = "  -3.06%  "
$ws.Range("E11").Value = "  -1.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.07"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.01%  "
$ws.Range("E13").Value = "  -1.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "676.96"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +11.37%  "
$ws.Range("D15").Value = "3.860.17"
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.46"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.92%  "
$ws.Range("D17").Value = "66.234.84"
$ws.Range("E17").Value = "  -0.33%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.90"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.29%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.117"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.63%  "
$ws.Range("D20").Value = "3.332.01"
$ws.Range("E20").Value = "  -0.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.09"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.896"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.79"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -5.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "102.03"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.05"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.94%  "
$ws.Range("E26").Value = "  -1.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.78"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("E28").Value = "  -3.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "32.31"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +5.67%  "
$ws.Range("E30").Value = "  -2.78%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.72"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "608.53"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +8.02%  "
$ws.Range("E33").Value = "  -2.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.12"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.52%  "
$ws.Range("E35").Value = "  -1.25%  "
$ws.Range("D36").Value = "3.823.04"
$ws.Range("E36").Value = "  +2.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.99"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.53%  "
$ws.Range("E39").Value = "  -3.16%  "
$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").Value = "0.0₃0702"
$ws.Range("E40").Value = "  -4.33%  "
$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.67"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.60%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.18"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "32.67"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -4.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.43"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +4.63%  "
$ws.Range("E45").Value = "  -2.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0415"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.99"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -12.67%  "
$ws.Range("E48").Value = "  -2.05%  "
$ws.Range("E49").Value = "  +0.30%  "
$ws.Range("E50").Value = "  -2.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "131.29"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +5.92%  "
